{"js": "// Update the worksheet date and all 100 arithmetic-answer cells in the\n// single 20x5 table to their new values (positional replacement \u2014 same\n// cell count/layout before and after the edit).\n\nconst newDate = \"2025-09-01 Monday\";\n\n// New values for the 20x5 answer table, row-major (row 0 .. row 19).\nconst newRows = [\n  [\"67-28=39\", \"35-34=1\", \"11+15=26\", \"26+46=72\", \"37-15=22\"],\n  [\"76-71=5\", \"45-11=34\", \"3+59=62\", \"1+73=74\", \"26+51=77\"],\n  [\"16+65=81\", \"68-3=65\", \"94-71=23\", \"60+13=73\", \"37-32=5\"],\n  [\"22+47=69\", \"12+43=55\", \"47+36=83\", \"26-0=26\", \"4+62=66\"],\n  [\"72-61=11\", \"53-36=17\", \"36+19=55\", \"32+46=78\", \"22+19=41\"],\n  [\"58+10=68\", \"71-15=56\", \"82+6=88\", \"1+97=98\", \"2+2=4\"],\n  [\"66-56=10\", \"56+33=89\", \"89-62=27\", \"66-6=60\", \"84-44=40\"],\n  [\"45+33=78\", \"99-30=69\", \"72-7=65\", \"39+44=83\", \"70+9=79\"],\n  [\"19+15=34\", \"59-15=44\", \"69-10=59\", \"47-38=9\", \"10+88=98\"],\n  [\"97-3=94\", \"54+33=87\", \"8-3=5\", \"34+54=88\", \"41+22=63\"],\n  [\"23+38=61\", \"18+4=22\", \"62-17=45\", \"63-40=23\", \"7+4=11\"],\n  [\"9+10=19\", \"29-6=23\", \"52+0=52\", \"75+10=85\", \"83+5=88\"],\n  [\"31-2=29\", \"49-13=36\", \"52+23=75\", \"90-15=75\", \"59-29=30\"],\n  [\"60-50=10\", \"76-6=70\", \"74-41=33\", \"46+4=50\", \"34-32=2\"],\n  [\"72-66=6\", \"78-71=7\", \"2+59=61\", \"30+65=95\", \"93-53=40\"],\n  [\"55-21=34\", \"57-46=11\", \"25+70=95\", \"18-11=7\", \"26+56=82\"],\n  [\"66-24=42\", \"47+20=67\", \"42+45=87\", \"42+57=99\", \"35+30=65\"],\n  [\"6+0=6\", \"4+52=56\", \"81-53=28\", \"24-8=16\", \"76-61=15\"],\n  [\"57-43=14\", \"12+75=87\", \"6+13=19\", \"73-2=71\", \"80+13=93\"],\n  [\"63-29=34\", \"57+4=61\", \"95-11=84\", \"15+9=24\", \"83-25=58\"],\n];\n\n// 1) Update the date paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text !== newDate) {\n  dateParagraph.getRange().insertText(newDate, Word.InsertLocation.replace);\n}\n\n// 2) Update every cell of the single answers table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < newRows[r].length; c++) {\n    const newValue = newRows[r][c];\n    if (table.values[r][c] !== newValue) {\n      const cell = table.getCell(r, c);\n      cell.value = newValue;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and all 100 arithmetic-answer cells in the\n# single 20x5 table to their new values (positional replacement \u2014 same\n# cell count/layout before and after the edit).\n\n$d = $word.ActiveDocument\n\n$newDate = \"2025-09-01 Monday\"\n\n# New values for the 20x5 answer table, row-major (row 1 .. row 20).\n$newRows = @(\n    @(\"67-28=39\", \"35-34=1\", \"11+15=26\", \"26+46=72\", \"37-15=22\"),\n    @(\"76-71=5\", \"45-11=34\", \"3+59=62\", \"1+73=74\", \"26+51=77\"),\n    @(\"16+65=81\", \"68-3=65\", \"94-71=23\", \"60+13=73\", \"37-32=5\"),\n    @(\"22+47=69\", \"12+43=55\", \"47+36=83\", \"26-0=26\", \"4+62=66\"),\n    @(\"72-61=11\", \"53-36=17\", \"36+19=55\", \"32+46=78\", \"22+19=41\"),\n    @(\"58+10=68\", \"71-15=56\", \"82+6=88\", \"1+97=98\", \"2+2=4\"),\n    @(\"66-56=10\", \"56+33=89\", \"89-62=27\", \"66-6=60\", \"84-44=40\"),\n    @(\"45+33=78\", \"99-30=69\", \"72-7=65\", \"39+44=83\", \"70+9=79\"),\n    @(\"19+15=34\", \"59-15=44\", \"69-10=59\", \"47-38=9\", \"10+88=98\"),\n    @(\"97-3=94\", \"54+33=87\", \"8-3=5\", \"34+54=88\", \"41+22=63\"),\n    @(\"23+38=61\", \"18+4=22\", \"62-17=45\", \"63-40=23\", \"7+4=11\"),\n    @(\"9+10=19\", \"29-6=23\", \"52+0=52\", \"75+10=85\", \"83+5=88\"),\n    @(\"31-2=29\", \"49-13=36\", \"52+23=75\", \"90-15=75\", \"59-29=30\"),\n    @(\"60-50=10\", \"76-6=70\", \"74-41=33\", \"46+4=50\", \"34-32=2\"),\n    @(\"72-66=6\", \"78-71=7\", \"2+59=61\", \"30+65=95\", \"93-53=40\"),\n    @(\"55-21=34\", \"57-46=11\", \"25+70=95\", \"18-11=7\", \"26+56=82\"),\n    @(\"66-24=42\", \"47+20=67\", \"42+45=87\", \"42+57=99\", \"35+30=65\"),\n    @(\"6+0=6\", \"4+52=56\", \"81-53=28\", \"24-8=16\", \"76-61=15\"),\n    @(\"57-43=14\", \"12+75=87\", \"6+13=19\", \"73-2=71\", \"80+13=93\"),\n    @(\"63-29=34\", \"57+4=61\", \"95-11=84\", \"15+9=24\", \"83-25=58\")\n)\n\n# 1) Update the date paragraph (first paragraph of the body).\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text -ne $newDate) {\n    $dateParagraph.Range.Text = $newDate\n}\n\n# 2) Update every cell of the single answers table.\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $newValue = $newRows[$r - 1][$c - 1]\n        $cell = $table.Cell($r, $c)\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $newValue) {\n            $cell.Range.Text = $newValue\n        }\n    }\n}\n"}
